# investigationOPU and cycle list pages changes
$wb = $excel.ActiveWorkbook

# --- HomePage: move selection ---
$wsHome = $wb.Worksheets.Item("HomePage")
$wsHome.Range("E20").Select()

# --- Investigation: move selection ---
$wsInv = $wb.Worksheets.Item("Investigation")
$wsInv.Range("M13").Select()

# --- Allergies: move selection ---
$wsAllerg = $wb.Worksheets.Item("Allergies")
$wsAllerg.Range("A13").Select()

# --- CycleList: rebuild the ART protocol dropdown list ---
$ws = $wb.Worksheets.Item("CycleList")

$ws.Range("D3").Value = "New Cycle"
$ws.Range("C1").Value = "ProtocolName"
$ws.Range("A1").Value = "ARTtype"
$ws.Range("B1").Value = "ARTSubtype"
$ws.Range("C2").Value = "Select"
$ws.Range("C4").Value = "Agonist Depot"
$ws.Range("C5").Value = "Antagonist"
$ws.Range("C6").Value = "Flare"
$ws.Range("C7").Value = "GnRh Long Protocol"
$ws.Range("C8").Value = "Minimal Stimulation"
$ws.Range("C9").Value = "Other"
$ws.Range("C10").Value = "Recipient Antagonist"
$ws.Range("C11").Value = "Recipient Depot"
$ws.Range("C12").Value = "Ultrashort"
$ws.Range("C13").Value = "Soft Protocol"
$ws.Range("C3").Value = "Agonist Daily"
$ws.Range("A2").Value = "OPU"
$ws.Range("D1").Value = "ListTitle"

# Match the bold/yellow-fill header style already used on the header row.
$ws.Range("A1:D1").Font.Bold = $true
$ws.Range("A1:D1").Interior.Color = 65535

$ws.Columns.Item(1).ColumnWidth = 8.78125
$ws.Columns.Item(2).ColumnWidth = 12
$ws.Columns.Item(3).ColumnWidth = 19.39453125
$ws.Columns.Item(4).ColumnWidth = 13.85546875

# Keep CycleList the active/selected tab, matching the source sheet.
$ws.Range("E10").Select()
